# Add a "notes" sheet (with the SB176 note) between "settings" and
# "file time fix", and leave it as the active/selected sheet - matching
# the author's commit "added note for SB176".

$wb = $excel.ActiveWorkbook

# Record the user's last selection on "settings" (C3) before switching away.
$settings = $wb.Worksheets.Item("settings")
$settings.Activate()
$settings.Range("C3").Select() | Out-Null

# Create the new sheet and populate it.
$notes = $wb.Worksheets.Add()
$notes.Name = "notes"
$notes.Range("A1").Value = "SB176 header needed to be corrected to match standard style"

# Reposition it right after "settings" (i.e. before "file time fix").
$notes.Move($null, $wb.Worksheets.Item("settings"))

# Re-fetch by name (Move invalidates the old reference / resets the active
# sheet), make it active, and select A2 underneath the note.
$notes = $wb.Worksheets.Item("notes")
$notes.Activate()
$notes.Range("A2").Select() | Out-Null
